$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "Complexity" (column C) values for the affected rows.
# Key = cell reference, Value = new numeric value to write.
$updates = @{
    "C134" = 1.160150658571203
    "C145" = 1.168885292845949
    "C146" = 1.180284271793053
    "C147" = 1.21486052238984
    "C148" = 1.276010857199003
    "C149" = 1.225847009048619
    "C150" = 1.223997632247695
    "C151" = 1.246643935962412
    "C152" = 1.236062985667228
    "C168" = 1.259233531562834
    "C169" = 1.257528759579208
    "C170" = 1.254173798833946
    "C171" = 1.251508208390669
    "C172" = 1.250056914663057
    "C173" = 1.26130329888185
    "C174" = 1.272806961088436
    "C175" = 1.261125837412918
    "C176" = 1.249623747082607
    "C177" = 1.249577153948315
    "C178" = 1.24948293342525
    "C179" = 1.24989081645374
    "C180" = 1.262322228572568
    "C181" = 1.240854830291638
    "C182" = 1.242066799159612
    "C183" = 1.24183318653802
    "C184" = 1.241279612420844
    "C185" = 1.240343896002161
    "C186" = 1.239565513799398
    "C188" = 1.22690209228363
    "C193" = 1.226995105162298
    "C195" = 1.202035027553715
    "C196" = 1.201453220846558
    "C197" = 1.192073308631313
    "C200" = 1.182850136424509
    "C201" = 1.134412140823284
    "C202" = 1.155892290501299
    "C206" = 1.160352133350743
    "C207" = 1.158309344632277
    "C209" = 1.198220063701776
    "C210" = 1.195588336915779
    "C214" = 1.219626455980141
    "C217" = 1.246098319938603
    "C218" = 1.246273242820492
    "C223" = 1.261691503158172
    "C225" = 1.209662791130737
    "C226" = 1.209156750441042
    "C231" = 1.209350538239505
    "C232" = 1.208327948655497
    "C233" = 1.207513369285311
    "C235" = 1.2541140304283
    "C236" = 1.252689620250149
    "C237" = 1.251961778778859
    "C238" = 1.251876098844296
    "C241" = 1.262663169287868
    "C243" = 1.299960470038234
    "C245" = 1.350402306315797
    "C246" = 1.374317018280406
    "C250" = 1.404521719344106
    "C251" = 1.488161909805113
    "C254" = 1.46652171772181
    "C257" = 1.529561818741367
    "C260" = 1.48817942768943
    "C261" = 1.500564367995939
    "C262" = 1.49977571503411
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}

Write-Host "Updated $($updates.Count) cells."
